$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-add the 'IDH1-M2-1-F' (mirrors the M column's header/data) as a new column O,
# copying both the header cell (text + style) and its data value.
$ws.Range("M1").Copy($ws.Range("O1"))
$ws.Range("M2").Copy($ws.Range("O2"))

# Restore the selection cursor position as recorded at save time
$ws.Range("L8").Select()
